# Bugfixed the naive forecaster component module
#
# Column A held Excel date serials (Dec-31 of each year) formatted with a
# custom "YYYY-MM-DD HH:MM:SS" number format. The forecaster's date axis
# should instead be labeled by fiscal-year quarter text, e.g. "1987Q4".
# Replace each date cell in A2:A39 with the equivalent "<year>Q4" text
# label and restyle it to match the header row (bold, centered, bordered,
# no date number format) instead of the old date-stamped style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startYear = 1987
$lastRow = 39

for ($row = 2; $row -le $lastRow; $row++) {
    $year = $startYear + ($row - 2)
    $ws.Cells.Item($row, 1).Value = "$($year)Q4"
}

# Re-use the header cell's existing style (bold font, border, centered)
# instead of the old date-format style, so no stray cell format is left
# behind on the quarter-label column.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
